# Apply updated LR-pair statistics (Natmi) per Dr Hou advice
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> column letter -> new value
$updates = @{}
$updates[2] = @{ 'E'=3; 'G'=46.70430066666666; 'H'=140.112902; 'I'=0.4277960227396158; 'J'=0.4350095176968582; 'K'=3; 'M'=10.68421466666667; 'N'=32.052644; 'O'=0.04705285980693976; 'P'=0.04892736897547583; 'Q'=498.9987741792098; 'R'=4490.988967612888; 'S'=0.02012902628393356; 'T'=0.02128387118019797 }
$updates[3] = @{ 'E'=3; 'G'=46.70430066666666; 'H'=140.112902; 'I'=0.4277960227396158; 'J'=0.4350095176968582; 'K'=3; 'M'=70.36235166666667; 'N'=211.087055; 'O'=0.3098730203341347; 'P'=0.3222178559101571; 'Q'=3286.224427853734; 'R'=29576.01985068361; 'S'=0.1325624456532549; 'T'=0.1401678340927932 }
$updates[4] = @{ 'E'=3; 'G'=46.70430066666666; 'H'=140.112902; 'I'=0.4277960227396158; 'J'=0.4350095176968582; 'K'=3; 'M'=59.09107466666666; 'N'=177.273224; 'O'=0.2602347611759026; 'P'=0.2706020894912812; 'Q'=2759.80731794845; 'R'=24838.26586153605; 'S'=0.1113273958096449; 'T'=0.1177144844373643 }
$updates[5] = @{ 'E'=3; 'G'=46.70430066666666; 'H'=140.112902; 'I'=0.4277960227396158; 'J'=0.4350095176968582; 'K'=3; 'M'=60.83231733333333; 'N'=182.496952; 'O'=0.2679031251727568; 'P'=0.2785759485989269; 'Q'=2841.130838986078; 'R'=25570.1775508747; 'S'=0.1146078914284188; 'T'=0.121183189041964 }
$updates[6] = @{ 'E'=3; 'G'=46.70430066666666; 'H'=140.112902; 'I'=0.4277960227396158; 'J'=0.4350095176968582; 'K'=2; 'M'=26.0983795; 'N'=52.196759; 'O'=0.1149362335102661; 'P'=0.07967673702415903; 'Q'=1218.90656308077; 'R'=7313.439378484618; 'S'=0.0491692635643636; 'T'=0.03466013894453883 }
$updates[7] = @{ 'E'=3; 'G'=14.05160533333333; 'H'=42.154816; 'I'=0.1287080801746603; 'J'=0.1308783553477452; 'K'=3; 'M'=10.68421466666667; 'N'=32.052644; 'O'=0.04705285980693976; 'P'=0.04892736897547583; 'Q'=150.1303677926116; 'R'=1351.173310133504; 'S'=0.006056083252478654; 'T'=0.006403533583002571 }
$updates[8] = @{ 'E'=3; 'G'=14.05160533333333; 'H'=42.154816; 'I'=0.1287080801746603; 'J'=0.1308783553477452; 'K'=3; 'M'=70.36235166666667; 'N'=211.087055; 'O'=0.3098730203341347; 'P'=0.3222178559101571; 'Q'=988.7039959452089; 'R'=8898.335963506879; 'S'=0.03988316154512995; 'T'=0.0421713430451981 }
$updates[9] = @{ 'E'=3; 'G'=14.05160533333333; 'H'=42.154816; 'I'=0.1287080801746603; 'J'=0.1308783553477452; 'K'=3; 'M'=59.09107466666666; 'N'=177.273224; 'O'=0.2602347611759026; 'P'=0.2706020894912812; 'Q'=830.3244599385315; 'R'=7472.920139446784; 'S'=0.03349431650566165; 'T'=0.03541595642628226 }
$updates[10] = @{ 'E'=3; 'G'=14.05160533333333; 'H'=42.154816; 'I'=0.1287080801746603; 'J'=0.1308783553477452; 'K'=3; 'M'=60.83231733333333; 'N'=182.496952; 'O'=0.2679031251727568; 'P'=0.2785759485989269; 'Q'=854.7917146800924; 'R'=7693.125432120831; 'S'=0.03448129691377724; 'T'=0.03645956199206556 }
$updates[11] = @{ 'E'=3; 'G'=14.05160533333333; 'H'=42.154816; 'I'=0.1287080801746603; 'J'=0.1308783553477452; 'K'=2; 'M'=26.0983795; 'N'=52.196759; 'O'=0.1149362335102661; 'P'=0.07967673702415903; 'Q'=366.7241285735573; 'R'=2200.344771441344; 'S'=0.01479322195761281; 'T'=0.01042796030119673 }
$updates[12] = @{ 'E'=3; 'G'=15.248849; 'H'=45.746547; 'I'=0.1396744381232708; 'J'=0.1420296279836289; 'K'=3; 'M'=10.68421466666667; 'N'=32.052644; 'O'=0.04705285980693976; 'P'=0.04892736897547583; 'Q'=162.9219761355853; 'R'=1466.297785220268; 'S'=0.006572081755627342; 'T'=0.006949136013804581 }
$updates[13] = @{ 'E'=3; 'G'=15.248849; 'H'=45.746547; 'I'=0.1396744381232708; 'J'=0.1420296279836289; 'K'=3; 'M'=70.36235166666667; 'N'=211.087055; 'O'=0.3098730203341347; 'P'=0.3222178559101571; 'Q'=1072.944875849898; 'R'=9656.503882649085; 'S'=0.04328134000473113; 'T'=0.04576448220460216 }
$updates[14] = @{ 'E'=3; 'G'=15.248849; 'H'=45.746547; 'I'=0.1396744381232708; 'J'=0.1420296279836289; 'K'=3; 'M'=59.09107466666666; 'N'=177.273224; 'O'=0.2602347611759026; 'P'=0.2706020894912812; 'Q'=901.0708748397253; 'R'=8109.637873557528; 'S'=0.03634814404738775; 'T'=0.03843351410203934 }
$updates[15] = @{ 'E'=3; 'G'=15.248849; 'H'=45.746547; 'I'=0.1396744381232708; 'J'=0.1420296279836289; 'K'=3; 'M'=60.83231733333333; 'N'=182.496952; 'O'=0.2679031251727568; 'P'=0.2785759485989269; 'Q'=927.6228213360826; 'R'=8348.605392024743; 'S'=0.03741921847997309; 'T'=0.03956603834469212 }
$updates[16] = @{ 'E'=3; 'G'=15.248849; 'H'=45.746547; 'I'=0.1396744381232708; 'J'=0.1420296279836289; 'K'=2; 'M'=26.0983795; 'N'=52.196759; 'O'=0.1149362335102661; 'P'=0.07967673702415903; 'Q'=397.9702481401955; 'R'=2387.821488841173; 'S'=0.01605365383555147; 'T'=0.01131645731849074 }
$updates[17] = @{ 'E'=3; 'G'=27.73836633333333; 'H'=83.215099; 'I'=0.2540743063339262; 'J'=0.258358506350017; 'K'=3; 'M'=10.68421466666667; 'N'=32.052644; 'O'=0.04705285980693976; 'P'=0.04892736897547583; 'Q'=296.3626604079729; 'R'=2667.263943671756; 'S'=0.0119549227164757; 'T'=0.0126408019681401 }
$updates[18] = @{ 'E'=3; 'G'=27.73836633333333; 'H'=83.215099; 'I'=0.2540743063339262; 'J'=0.258358506350017; 'K'=3; 'M'=70.36235166666667; 'N'=211.087055; 'O'=0.3098730203341347; 'P'=0.3222178559101571; 'Q'=1951.736686604827; 'R'=17565.63017944344; 'S'=0.07873077269299389; 'T'=0.08324772397225318 }
$updates[19] = @{ 'E'=3; 'G'=27.73836633333333; 'H'=83.215099; 'I'=0.2540743063339262; 'J'=0.258358506350017; 'K'=3; 'M'=59.09107466666666; 'N'=177.273224; 'O'=0.2602347611759026; 'P'=0.2706020894912812; 'Q'=1639.089876134353; 'R'=14751.80888520918; 'S'=0.0661189664297424; 'T'=0.06991235165616105 }
$updates[20] = @{ 'E'=3; 'G'=27.73836633333333; 'H'=83.215099; 'I'=0.2540743063339262; 'J'=0.258358506350017; 'K'=3; 'M'=60.83231733333333; 'N'=182.496952; 'O'=0.2679031251727568; 'P'=0.2785759485989269; 'Q'=1687.389103097583; 'R'=15186.50192787825; 'S'=0.06806730069295919; 'T'=0.07197246598505785 }
$updates[21] = @{ 'E'=3; 'G'=27.73836633333333; 'H'=83.215099; 'I'=0.2540743063339262; 'J'=0.258358506350017; 'K'=2; 'M'=26.0983795; 'N'=52.196759; 'O'=0.1149362335102661; 'P'=0.07967673702415903; 'Q'=723.9264112773568; 'R'=4343.558467664141; 'S'=0.02920234380175503; 'T'=0.02058516276840483 }
$updates[22] = @{ 'E'=2; 'G'=5.431107000000001; 'H'=10.862214; 'I'=0.0497471526285271; 'J'=0.03372399262175058; 'K'=3; 'M'=10.68421466666667; 'N'=32.052644; 'O'=0.04705285980693976; 'P'=0.04892736897547583; 'Q'=58.02711306563602; 'R'=348.1626783938161; 'S'=0.00234074579842452; 'T'=0.001650026230330615 }
$updates[23] = @{ 'E'=2; 'G'=5.431107000000001; 'H'=10.862214; 'I'=0.0497471526285271; 'J'=0.03372399262175058; 'K'=3; 'M'=70.36235166666667; 'N'=211.087055; 'O'=0.3098730203341347; 'P'=0.3222178559101571; 'Q'=382.1454606732951; 'R'=2292.87276403977; 'S'=0.01541530043802488; 'T'=0.01086647259531043 }
$updates[24] = @{ 'E'=2; 'G'=5.431107000000001; 'H'=10.862214; 'I'=0.0497471526285271; 'J'=0.03372399262175058; 'K'=3; 'M'=59.09107466666666; 'N'=177.273224; 'O'=0.2602347611759026; 'P'=0.2706020894912812; 'Q'=320.929949259656; 'R'=1925.579695557936; 'S'=0.01294593838346593; 'T'=0.00912578286943426 }
$updates[25] = @{ 'E'=2; 'G'=5.431107000000001; 'H'=10.862214; 'I'=0.0497471526285271; 'J'=0.03372399262175058; 'K'=3; 'M'=60.83231733333333; 'N'=182.496952; 'O'=0.2679031251727568; 'P'=0.2785759485989269; 'Q'=330.386824495288; 'R'=1982.320946971728; 'S'=0.01332741765762853; 'T'=0.00939469323514738 }
$updates[26] = @{ 'E'=2; 'G'=5.431107000000001; 'H'=10.862214; 'I'=0.0497471526285271; 'J'=0.03372399262175058; 'K'=2; 'M'=26.0983795; 'N'=52.196759; 'O'=0.1149362335102661; 'P'=0.07967673702415903; 'Q'=141.7430915911065; 'R'=566.9723663644261; 'S'=0.005717750350983241; 'T'=0.002687017691527901 }

foreach ($rowNum in $updates.Keys) {
    $rowUpdates = $updates[$rowNum]
    foreach ($col in $rowUpdates.Keys) {
        $cellRef = "$col$rowNum"
        $ws.Range($cellRef).Value = $rowUpdates[$col]
    }
}

Write-Host "Applied $($updates.Keys.Count) row updates across $($updates.Keys | ForEach-Object { $updates[$_].Keys.Count } | Measure-Object -Sum | Select-Object -ExpandProperty Sum) cells"